$wb = $excel.ActiveWorkbook

# --- awards sheet: insert a new row 2 for the 2020 award ---
$ws3 = $wb.Worksheets.Item("awards")
$ws3.Rows("2:2").Insert()
$ws3.Range("A2").Value = "2020-09-17"
$ws3.Range("B2").Value = "2020"
$ws3.Range("D2").Value = "Gaetano Borriello Outstanding Award Finalist"
$ws3.Range("C2").Value = "UbiComp/ISWC 2020"
$ws3.Range("E2").Value = "Designing Interactive Technologies to Encourage Physical Activities for Health Behavior Promotion"

# --- academicServices sheet: append new row 5 ---
$ws1 = $wb.Worksheets.Item("academicServices")
$ws1.Range("A4").Copy($ws1.Range("A5"))
$ws1.Range("A5").Value = 43983
$ws1.Range("B5").Value = 2020
$ws1.Range("C5").Value = "ACM HEALTH"
$ws1.Range("D4").Copy($ws1.Range("D5"))
$ws1.Range("D5").Value = "External reviewer"

# --- restore/update the on-screen selections to match the edited session ---
[void]$ws3.Range("E11").Select()
[void]$ws1.Range("D8").Select()
